# Apply cell value updates per the target diff (odds/line movements).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.75
$ws.Range("H2").Value = 4.2
$ws.Range("I2").Value = 4.5
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 1.19
$ws.Range("R2").Value = 1.97
$ws.Range("Y2").Value = 36
$ws.Range("Z2").Value = 50
$ws.Range("AA2").Value = 100
$ws.Range("AC2").Value = 15
$ws.Range("AD2").Value = 980
$ws.Range("AE2").Value = 48
$ws.Range("AF2").Value = 20
$ws.Range("AH2").Value = 1000
$ws.Range("AI2").Value = 42
$ws.Range("AJ2").Value = 22
$ws.Range("AK2").Value = 18.5
$ws.Range("AL2").Value = 26
$ws.Range("AN2").Value = 5.6
# Row 3
$ws.Range("F3").Value = 2.78
$ws.Range("I3").Value = 3.5
$ws.Range("J3").Value = 2.4
$ws.Range("L3").Value = 1.58
$ws.Range("V3").Value = 1.4
$ws.Range("Z3").Value = 21
$ws.Range("AA3").Value = 70
$ws.Range("AE3").Value = 65
$ws.Range("AF3").Value = 25
$ws.Range("AH3").Value = 36
$ws.Range("AJ3").Value = 90
$ws.Range("AK3").Value = 80
# Row 4
$ws.Range("G4").Value = 6.2
$ws.Range("I4").Value = 1.95
$ws.Range("J4").Value = 3.4
$ws.Range("L4").Value = 1.41
$ws.Range("M4").Value = 1.08
$ws.Range("R4").Value = 1.26
$ws.Range("U4").Value = 1.83
$ws.Range("V4").Value = 2.04
$ws.Range("Y4").Value = 8.4
$ws.Range("AC4").Value = 970
$ws.Range("AI4").Value = 980
$ws.Range("AK4").Value = 100
$ws.Range("AN4").Value = 140
# Row 5
$ws.Range("K5").Value = 6.2
$ws.Range("Q5").Value = 1.55
$ws.Range("T5").Value = 1.94
$ws.Range("U5").Value = 1.95
$ws.Range("AC5").Value = 14.5
$ws.Range("AL5").Value = 140
$ws.Range("AN5").Value = 200
# Row 6
$ws.Range("F6").Value = 2.66
$ws.Range("G6").Value = 2.76
$ws.Range("H6").Value = 3.05
$ws.Range("I6").Value = 3.15
$ws.Range("J6").Value = 3.2
$ws.Range("K6").Value = 3.3
$ws.Range("P6").Value = 1.69
$ws.Range("Q6").Value = 2.38
$ws.Range("T6").Value = 1.97
$ws.Range("U6").Value = 1.97
$ws.Range("AA6").Value = 70
$ws.Range("AB6").Value = 9
$ws.Range("AI6").Value = 75
# Row 7
$ws.Range("J7").Value = 8.4
$ws.Range("Q7").Value = 1.59
$ws.Range("R7").Value = 1.63
$ws.Range("AM7").Value = 440
# Row 8
$ws.Range("F8").Value = 2.58
$ws.Range("H8").Value = 2.98
$ws.Range("J8").Value = 3.35
$ws.Range("Y8").Value = 11.5
$ws.Range("AA8").Value = 980
$ws.Range("AG8").Value = 12.5
$ws.Range("AK8").Value = 30
$ws.Range("AO8").Value = 38
# Row 9
$ws.Range("G9").Value = 2.66
$ws.Range("H9").Value = 3.7
$ws.Range("K9").Value = 3
# Row 10
$ws.Range("I10").Value = 4.4
